# Update cryptos list with latest prices and 1h volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.418.66"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.584.14"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'213.40"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'0.492"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'44.64"
$ws.Range("D9").Value = "'23.92"
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "1.582.59"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").Value = "'0.518"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "28.458.89"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "'62.14"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").Value = "'230.05"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "'7.45"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  -2.42%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'3.91"
$ws.Range("E23").Value = "  -3.30%  "
$ws.Range("D24").Value = "'9.15"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("D26").Value = "'151.67"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "'15.03"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'0.0482"
$ws.Range("E31").Value = "  +2.56%  "
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("D35").Value = "1.394.38"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("E36").Value = "  +6.71%  "
$ws.Range("E37").Value = "  -4.89%  "
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").Value = "'0.522"
$ws.Range("E41").Value = "  -3.40%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'0.791"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").Value = "'5.44"
$ws.Range("E45").Value = "  -3.37%  "
$ws.Range("D46").Value = "'0.0457"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("D47").Value = "'0.959"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("D48").Value = "'62.92"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "1.723.51"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'86.61"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "0.0₆0102"
$ws.Range("E51").Value = "  -2.52%  "
